$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi following Dr Hou advice: refresh LR-pair values for Tgfb2-Tgfbr3,
# adding the "M2" target cluster and updating all computed NATMI metrics
# for the full 3 (sending) x 4 (target) cluster combinations.

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tgfb2"
$ws.Range("C2").Value = "Tgfbr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.191602
$ws.Range("H2").Value = 6.574806000000001
$ws.Range("I2").Value = 0.07674610985252207
$ws.Range("J2").Value = 0.07674610985252209
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 23.68145933333333
$ws.Range("N2").Value = 71.04437799999999
$ws.Range("O2").Value = 0.2340556429375698
$ws.Range("P2").Value = 0.2340556429375698
$ws.Range("Q2").Value = 51.900333637852
$ws.Range("R2").Value = 467.103002740668
$ws.Range("S2").Value = 0.01796286008448942
$ws.Range("T2").Value = 0.01796286008448941

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tgfb2"
$ws.Range("C3").Value = "Tgfbr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.191602
$ws.Range("H3").Value = 6.574806000000001
$ws.Range("I3").Value = 0.07674610985252207
$ws.Range("J3").Value = 0.07674610985252209
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 62.46631733333334
$ws.Range("N3").Value = 187.398952
$ws.Range("O3").Value = 0.6173856880862099
$ws.Range("P3").Value = 0.6173856880862097
$ws.Range("Q3").Value = 136.901306000368
$ws.Range("R3").Value = 1232.111754003312
$ws.Range("S3").Value = 0.04738194983923919
$ws.Range("T3").Value = 0.04738194983923919

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tgfb2"
$ws.Range("C4").Value = "Tgfbr3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.191602
$ws.Range("H4").Value = 6.574806000000001
$ws.Range("I4").Value = 0.07674610985252207
$ws.Range("J4").Value = 0.07674610985252209
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.08293366666666667
$ws.Range("N4").Value = 0.248801
$ws.Range("O4").Value = 0.0008196746830341778
$ws.Range("P4").Value = 0.0008196746830341776
$ws.Range("Q4").Value = 0.181757589734
$ws.Range("R4").Value = 1.635818307606
$ws.Range("S4").Value = 0.00006290684326747222
$ws.Range("T4").Value = 0.00006290684326747222

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Tgfb2"
$ws.Range("C5").Value = "Tgfbr3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.191602
$ws.Range("H5").Value = 6.574806000000001
$ws.Range("I5").Value = 0.07674610985252207
$ws.Range("J5").Value = 0.07674610985252209
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 14.948048
$ws.Range("N5").Value = 44.844144
$ws.Range("O5").Value = 0.1477389942931862
$ws.Range("P5").Value = 0.1477389942931862
$ws.Range("Q5").Value = 32.760171892896
$ws.Range("R5").Value = 294.841547036064
$ws.Range("S5").Value = 0.011338393085526
$ws.Range("T5").Value = 0.011338393085526

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tgfb2"
$ws.Range("C6").Value = "Tgfbr3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 16.59481266666667
$ws.Range("H6").Value = 49.78443799999999
$ws.Range("I6").Value = 0.5811216251390648
$ws.Range("J6").Value = 0.5811216251390647
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 23.68145933333333
$ws.Range("N6").Value = 71.04437799999999
$ws.Range("O6").Value = 0.2340556429375698
$ws.Range("P6").Value = 0.2340556429375698
$ws.Range("Q6").Value = 392.9893813099515
$ws.Range("R6").Value = 3536.904431789563
$ws.Range("S6").Value = 0.1360147955968493
$ws.Range("T6").Value = 0.1360147955968492

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tgfb2"
$ws.Range("C7").Value = "Tgfbr3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 16.59481266666667
$ws.Range("H7").Value = 49.78443799999999
$ws.Range("I7").Value = 0.5811216251390648
$ws.Range("J7").Value = 0.5811216251390647
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 62.46631733333334
$ws.Range("N7").Value = 187.398952
$ws.Range("O7").Value = 0.6173856880862099
$ws.Range("P7").Value = 0.6173856880862097
$ws.Range("Q7").Value = 1036.61683412322
$ws.Range("R7").Value = 9329.551507108976
$ws.Range("S7").Value = 0.358776174398258
$ws.Range("T7").Value = 0.3587761743982579

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Tgfb2"
$ws.Range("C8").Value = "Tgfbr3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 16.59481266666667
$ws.Range("H8").Value = 49.78443799999999
$ws.Range("I8").Value = 0.5811216251390648
$ws.Range("J8").Value = 0.5811216251390647
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.08293366666666667
$ws.Range("N8").Value = 0.248801
$ws.Range("O8").Value = 0.0008196746830341778
$ws.Range("P8").Value = 0.0008196746830341776
$ws.Range("Q8").Value = 1.376268662093111
$ws.Range("R8").Value = 12.386417958838
$ws.Range("S8").Value = 0.0004763306838901692
$ws.Range("T8").Value = 0.000476330683890169

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Tgfb2"
$ws.Range("C9").Value = "Tgfbr3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 16.59481266666667
$ws.Range("H9").Value = 49.78443799999999
$ws.Range("I9").Value = 0.5811216251390648
$ws.Range("J9").Value = 0.5811216251390647
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 14.948048
$ws.Range("N9").Value = 44.844144
$ws.Range("O9").Value = 0.1477389942931862
$ws.Range("P9").Value = 0.1477389942931862
$ws.Range("Q9").Value = 248.0600562923413
$ws.Range("R9").Value = 2232.540506631072
$ws.Range("S9").Value = 0.0858543244600674
$ws.Range("T9").Value = 0.08585432446006737

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tgfb2"
$ws.Range("C10").Value = "Tgfbr3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 9.770107666666666
$ws.Range("H10").Value = 29.310323
$ws.Range("I10").Value = 0.342132265008413
$ws.Range("J10").Value = 0.342132265008413
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 23.68145933333333
$ws.Range("N10").Value = 71.04437799999999
$ws.Range("O10").Value = 0.2340556429375698
$ws.Range("P10").Value = 0.2340556429375698
$ws.Range("Q10").Value = 231.3704073904549
$ws.Range("R10").Value = 2082.333666514094
$ws.Range("S10").Value = 0.08007798725623114
$ws.Range("T10").Value = 0.08007798725623111

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Tgfb2"
$ws.Range("C11").Value = "Tgfbr3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 9.770107666666666
$ws.Range("H11").Value = 29.310323
$ws.Range("I11").Value = 0.342132265008413
$ws.Range("J11").Value = 0.342132265008413
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 62.46631733333334
$ws.Range("N11").Value = 187.398952
$ws.Range("O11").Value = 0.6173856880862099
$ws.Range("P11").Value = 0.6173856880862097
$ws.Range("Q11").Value = 610.3026458868329
$ws.Range("R11").Value = 5492.723812981496
$ws.Range("S11").Value = 0.2112275638487126
$ws.Range("T11").Value = 0.2112275638487126

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Tgfb2"
$ws.Range("C12").Value = "Tgfbr3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 9.770107666666666
$ws.Range("H12").Value = 29.310323
$ws.Range("I12").Value = 0.342132265008413
$ws.Range("J12").Value = 0.342132265008413
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.08293366666666667
$ws.Range("N12").Value = 0.248801
$ws.Range("O12").Value = 0.0008196746830341778
$ws.Range("P12").Value = 0.0008196746830341776
$ws.Range("Q12").Value = 0.8102708525247778
$ws.Range("R12").Value = 7.292437672722999
$ws.Range("S12").Value = 0.0002804371558765363
$ws.Range("T12").Value = 0.0002804371558765362

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Tgfb2"
$ws.Range("C13").Value = "Tgfbr3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 9.770107666666666
$ws.Range("H13").Value = 29.310323
$ws.Range("I13").Value = 0.342132265008413
$ws.Range("J13").Value = 0.342132265008413
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 14.948048
$ws.Range("N13").Value = 44.844144
$ws.Range("O13").Value = 0.1477389942931862
$ws.Range("P13").Value = 0.1477389942931862
$ws.Range("Q13").Value = 146.0440383665013
$ws.Range("R13").Value = 1314.396345298512
$ws.Range("S13").Value = 0.05054627674759281
$ws.Range("T13").Value = 0.05054627674759279

